$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.839.20"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "3.508.20"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.52"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.89"
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("D7").Value = "3.506.72"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.10"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.383"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "4.118.53"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.16"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "3.521.63"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.116"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "64.955.29"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.67"
$ws.Range("E19").Value = "  -3.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.37"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.69"
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "389.64"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.575"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.65"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.656.44"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000112"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.78"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.60"
$ws.Range("E29").Value = "  +15.76%  "
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.28"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.37"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").Value = "3.518.74"
$ws.Range("E33").Value = "  -0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.99"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.144"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.24"
$ws.Range("E37").Value = "  +4.92%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "169.84"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.78"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0822"
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.818"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.74"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.57"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.23"
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.40"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.65"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.88"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").Value = "2.369.49"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0266"
$ws.Range("E51").Value = "  +2.31%  "
